# Insert a new price-report row for Feria Lagunitas de Puerto Montt - Mango
# at row 64 (weekly update), pushing the existing rows 64-151 down to 65-152.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 64; everything currently at/after row 64
# (including the old row 64) shifts down by one, growing the sheet from
# A1:T151 to A1:T152.
$ws.Rows(64).Insert()

# Populate the newly inserted row 64 with the new weekly observation.
$ws.Range("A64").Value = 4
$ws.Range("B64").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C64").Value = "Los Lagos"
$ws.Range("D64").Value = 44579
$ws.Range("E64").Value = 10
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100108
$ws.Range("H64").Value = "Tropicales y subtropicales"
$ws.Range("I64").Value = 100108002
$ws.Range("J64").Value = "Mango"
$ws.Range("K64").Value = "Sin especificar"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 180
$ws.Range("N64").Value = 8000
$ws.Range("O64").Value = 8500
$ws.Range("P64").Value = 8250
$ws.Range("Q64").Value = "$/bandeja 4 kilos"
$ws.Range("R64").Value = "Perú"
$ws.Range("S64").Value = 2062
$ws.Range("T64").Value = 4
